$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hotel_info")

# Update English_Reviews_num (G2), Local_Rank (H2) and Total_Reviews_num (I2)
# for the hotel row. These values are stored as text in the sheet, so force
# a text number format before assigning to avoid them being interpreted as
# numeric values.
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "3"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "9"

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "3"
